$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 50: new log entry - "Menu verbessert und kleine Verschönerungen im Code",
# 3 hours, dated "6-7-8-Nov"
$ws.Range("F50").Value = "6-7-8-Nov"
$ws.Range("B50").Value = "Menu verbessert und kleine Verschönerungen im Code"
$ws.Range("D50").Value = 3

# Row 52: new log entry - "Erweiterung des Kollisionssystem mit Swept AABB",
# 6 hours, dated 2018-11-09 (serial 43413), formatted like the other dates
$ws.Range("B52").Value = "Erweiterung des Kollisionssystem mit Swept AABB"
$ws.Range("D52").Value = 6
$ws.Range("F52").Value = 43413
$ws.Range("F52").NumberFormat = $ws.Range("F48").NumberFormat

# Sheet scrolled down and new selection reflecting the added rows
$ws.Application.ActiveWindow.ScrollRow = 37
$ws.Range("B53").Select()
